# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update / replace the worker data rows (B15:J25 table) with the new
# set of account-statement rows: each worker now appears for periods
# 1608 then 1607 (grouped by worker instead of by period).
$ws.Range("C16").Value = "9145960"
$ws.Range("D16").Value = "YOBANIS LORENZO MARTINEZ ZARZA"
$ws.Range("E16").Value = "1608"

$ws.Range("C17").Value = "9145960"
$ws.Range("D17").Value = "YOBANIS LORENZO MARTINEZ ZARZA"
$ws.Range("E17").Value = "1607"

$ws.Range("C18").Value = "9176750"
$ws.Range("D18").Value = "FELIPE SEGUNDO PADILLA DIAZ"
$ws.Range("E18").Value = "1608"

$ws.Range("C19").Value = "9176750"
$ws.Range("D19").Value = "FELIPE SEGUNDO PADILLA DIAZ"
$ws.Range("E19").Value = "1607"

$ws.Range("C20").Value = "45489734"
$ws.Range("D20").Value = "ASTRID GUZMAN POMPEYO"
$ws.Range("E20").Value = "1608"

$ws.Range("C21").Value = "45489734"
$ws.Range("D21").Value = "ASTRID GUZMAN POMPEYO"
$ws.Range("E21").Value = "1607"

$ws.Range("C22").Value = "22810873"
$ws.Range("D22").Value = "MAIRA PATRICIA PADILLA CABARCAS"
$ws.Range("E22").Value = "1608"

$ws.Range("C23").Value = "22810873"
$ws.Range("D23").Value = "MAIRA PATRICIA PADILLA CABARCAS"
$ws.Range("E23").Value = "1607"

$ws.Range("C24").Value = "1047478843"
$ws.Range("D24").Value = "JAVIER EDUARDO GONZALEZ POLO"
$ws.Range("E24").Value = "1608"

$ws.Range("C25").Value = "1047478843"
$ws.Range("D25").Value = "JAVIER EDUARDO GONZALEZ POLO"
$ws.Range("E25").Value = "1607"

$wb.Save()
